# "Tabla actualizada al español"
# Replace the placeholder "08:05AM" text/time values in column H (horaEntrada)
# with the real, differentiated entry-time values for each row, and update
# the sheet's active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2 already stored a real time value (08:05 AM) - correct it to 07:00 AM.
$ws.Range("H2").Value = 0.29166666666666669

# H3:H17 previously held the literal text "08:05AM" (a shared string used as a
# placeholder). Replace each with its real numeric time-of-day fraction so the
# cells become proper Excel times (keeping the existing hh:mm AM/PM format).
$ws.Range("H3").Value  = 0.33333333333333331
$ws.Range("H4").Value  = 0.33333333333333331
$ws.Range("H5").Value  = 0.58333333333333337
$ws.Range("H6").Value  = 0.041666666666666664
$ws.Range("H7").Value  = 0.041666666666666664
$ws.Range("H8").Value  = 0.25
$ws.Range("H9").Value  = 0.375
$ws.Range("H10").Value = 0.375
$ws.Range("H11").Value = 0.375
$ws.Range("H12").Value = 0.29166666666666669
$ws.Range("H13").Value = 0.29166666666666669
$ws.Range("H14").Value = 0.79166666666666663
$ws.Range("H15").Value = 0.33333333333333331
$ws.Range("H16").Value = 0.33333333333333331
$ws.Range("H17").Value = 0.041666666666666664

# Update the sheet's active cell / selection and scroll position to match
# the saved view state (topLeftCell moves from G1 to N1, selection to F9).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1
$null = $ws.Range("F9").Select()
